$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H5").Value = 345.5
$ws.Range("I5").Value = 345.5
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 345.5
$ws.Range("L5").Value = 0
$ws.Range("M5").Value = -230.5
$ws.Range("H37").Value = 500
$ws.Range("J37").Value = 500
$ws.Range("L37").Value = 1500
$ws.Range("N37").Value = -1752
$ws.Range("H48").Value = 2604.75
$ws.Range("I48").Value = 1800
$ws.Range("J48").Value = 2873
$ws.Range("K48").Value = 5400
$ws.Range("L48").Value = 8619
$ws.Range("M48").Value = -5108
$ws.Range("N48").Value = -9203
$ws.Range("H56").Value = 2604.75
$ws.Range("I56").Value = 1800
$ws.Range("J56").Value = 2873
$ws.Range("K56").Value = 5400
$ws.Range("L56").Value = 8619
$ws.Range("M56").Value = -4866
$ws.Range("N56").Value = -9687
$ws.Range("H107").Value = 737.19354
$ws.Range("I107").Value = 702.7083
$ws.Range("J107").Value = 855.4286
$ws.Range("K107").Value = 702.7083
$ws.Range("L107").Value = 855.4286
$ws.Range("M107").Value = 1217.2917
$ws.Range("N107").Value = -4695.4286
$ws.Range("H116").Value = 3939.95
$ws.Range("I116").Value = 3588
$ws.Range("J116").Value = 4291.9
$ws.Range("K116").Value = 3588
$ws.Range("L116").Value = 4291.9
$ws.Range("M116").Value = -146
$ws.Range("N116").Value = -11175.9
$ws.Range("H132").Value = 6670627.5
$ws.Range("I132").Value = 10003496
$ws.Range("J132").Value = 4890
$ws.Range("K132").Value = 30010488
$ws.Range("L132").Value = 14670
$ws.Range("M132").Value = -30007958
$ws.Range("N132").Value = -19730
$ws.Range("H137").Value = 3584.261
$ws.Range("I137").Value = 3796.879
$ws.Range("K137").Value = 11390.637
$ws.Range("M137").Value = -8840.636999999999
$ws.Range("H141").Value = 656250.7
$ws.Range("I141").Value = 1940.7142
$ws.Range("J141").Value = 1114267.8
$ws.Range("K141").Value = 5822.142599999999
$ws.Range("L141").Value = 3342803.4
$ws.Range("M141").Value = -642.1425999999992
$ws.Range("N141").Value = -3353163.4
$ws.Range("N5").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 5683224
$ws.Range("I2").Value = 9616472
$ws.Range("J2").Value = 1866.6666
$ws.Range("K2").Value = 9616472
$ws.Range("L2").Value = 1866.6666
$ws.Range("M2").Value = -9616359
$ws.Range("N2").Value = -2092.6666
$ws.Range("H45").Value = 1586.5667
$ws.Range("I45").Value = 1235.6072
$ws.Range("J45").Value = 6500
$ws.Range("K45").Value = 1235.6072
$ws.Range("L45").Value = 6500
$ws.Range("M45").Value = -858.6071999999999
$ws.Range("N45").Value = -7254
$ws.Range("H61").Value = 1953.6086
$ws.Range("I61").Value = 1495.2903
$ws.Range("K61").Value = 1495.2903
$ws.Range("M61").Value = -1283.2903
$ws.Range("H63").Value = 2377.5
$ws.Range("I63").Value = 1271.0714
$ws.Range("J63").Value = 6250
$ws.Range("K63").Value = 1271.0714
$ws.Range("L63").Value = 6250
$ws.Range("M63").Value = -585.0714
$ws.Range("N63").Value = -7622
$ws.Range("H66").Value = 2377.5
$ws.Range("I66").Value = 1271.0714
$ws.Range("J66").Value = 6250
$ws.Range("K66").Value = 6355.357
$ws.Range("L66").Value = 31250
$ws.Range("M66").Value = -2923.357
$ws.Range("N66").Value = -38114
$ws.Range("H102").Value = 2415.8333
$ws.Range("I102").Value = 2399
$ws.Range("K102").Value = 2399
$ws.Range("M102").Value = -777
$ws.Range("H110").Value = 2198.4443
$ws.Range("I110").Value = 672.36365
$ws.Range("J110").Value = 4596.5713
$ws.Range("K110").Value = 672.36365
$ws.Range("L110").Value = 4596.5713
$ws.Range("M110").Value = 1372.63635
$ws.Range("N110").Value = -8686.5713
$ws.Range("H116").Value = 5683224
$ws.Range("I116").Value = 9616472
$ws.Range("J116").Value = 1866.6666
$ws.Range("K116").Value = 9616472
$ws.Range("L116").Value = 1866.6666
$ws.Range("M116").Value = -9614178
$ws.Range("N116").Value = -6454.6666
$ws.Range("H122").Value = 1892.8148
$ws.Range("I122").Value = 1587.7368
$ws.Range("J122").Value = 2617.375
$ws.Range("K122").Value = 4763.2104
$ws.Range("L122").Value = 7852.125
$ws.Range("M122").Value = -2313.2104
$ws.Range("N122").Value = -12752.125
$ws.Range("H132").Value = 2564.261
$ws.Range("I132").Value = 1911.0857
$ws.Range("J132").Value = 4642.5454
$ws.Range("K132").Value = 5733.257100000001
$ws.Range("L132").Value = 13927.6362
$ws.Range("M132").Value = -3203.257100000001
$ws.Range("N132").Value = -18987.6362
$ws.Range("H136").Value = 1953.6086
$ws.Range("I136").Value = 1495.2903
$ws.Range("K136").Value = 4485.8709
$ws.Range("M136").Value = -1935.8709

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 5683224
$ws.Range("I3").Value = 9616472
$ws.Range("J3").Value = 1866.6666
$ws.Range("K3").Value = 9616472
$ws.Range("L3").Value = 1866.6666
$ws.Range("M3").Value = -9616358
$ws.Range("N3").Value = -2094.6666
$ws.Range("H99").Value = 2825.1177
$ws.Range("I99").Value = 2926.6924
$ws.Range("J99").Value = 2495
$ws.Range("K99").Value = 2926.6924
$ws.Range("L99").Value = 2495
$ws.Range("M99").Value = -1428.6924
$ws.Range("N99").Value = -5491
$ws.Range("H134").Value = 2716.3147
$ws.Range("I134").Value = 2801.9211
$ws.Range("K134").Value = 8405.763300000001
$ws.Range("M134").Value = -5870.763300000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 2206.8235
$ws.Range("I16").Value = 1631.6
$ws.Range("J16").Value = 3028.5715
$ws.Range("K16").Value = 1631.6
$ws.Range("L16").Value = 3028.5715
$ws.Range("M16").Value = -1344.6
$ws.Range("N16").Value = -3602.5715
$ws.Range("H31").Value = 2081.4036
$ws.Range("I31").Value = 1556.9231
$ws.Range("J31").Value = 7536
$ws.Range("K31").Value = 1556.9231
$ws.Range("L31").Value = 7536
$ws.Range("M31").Value = -1261.9231
$ws.Range("N31").Value = -8126
$ws.Range("H34").Value = 2081.4036
$ws.Range("I34").Value = 1556.9231
$ws.Range("J34").Value = 7536
$ws.Range("K34").Value = 1556.9231
$ws.Range("L34").Value = 7536
$ws.Range("M34").Value = -1354.9231
$ws.Range("N34").Value = -7940
$ws.Range("H99").Value = 2391.6667
$ws.Range("I99").Value = 1587.4286
$ws.Range("J99").Value = 3517.6
$ws.Range("K99").Value = 1587.4286
$ws.Range("L99").Value = 3517.6
$ws.Range("M99").Value = -89.42859999999996
$ws.Range("N99").Value = -6513.6
$ws.Range("H107").Value = 1179.44
$ws.Range("I107").Value = 940.375
$ws.Range("J107").Value = 1604.4445
$ws.Range("K107").Value = 940.375
$ws.Range("L107").Value = 1604.4445
$ws.Range("M107").Value = 979.625
$ws.Range("N107").Value = -5444.4445
$ws.Range("H113").Value = 2206.8235
$ws.Range("I113").Value = 1631.6
$ws.Range("J113").Value = 3028.5715
$ws.Range("K113").Value = 1631.6
$ws.Range("L113").Value = 3028.5715
$ws.Range("M113").Value = 538.4000000000001
$ws.Range("N113").Value = -7368.5715
$ws.Range("H122").Value = 2507
$ws.Range("I122").Value = 2071
$ws.Range("K122").Value = 6213
$ws.Range("M122").Value = -3763
$ws.Range("H126").Value = 2391.6667
$ws.Range("I126").Value = 1587.4286
$ws.Range("J126").Value = 3517.6
$ws.Range("K126").Value = 4762.2858
$ws.Range("L126").Value = 10552.8
$ws.Range("M126").Value = -2292.2858
$ws.Range("N126").Value = -15492.8
$ws.Range("H132").Value = 2265
$ws.Range("I132").Value = 2063.9048
$ws.Range("K132").Value = 6191.714399999999
$ws.Range("M132").Value = -3661.714399999999
$ws.Range("H134").Value = 14709776
$ws.Range("I134").Value = 20837682
$ws.Range("K134").Value = 62513046
$ws.Range("M134").Value = -62510511

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 4465.3335
$ws.Range("I70").Value = 4338
$ws.Range("J70").Value = 4720
$ws.Range("K70").Value = 4338
$ws.Range("L70").Value = 4720
$ws.Range("M70").Value = -4068
$ws.Range("N70").Value = -5260
$ws.Range("H73").Value = 4465.3335
$ws.Range("I73").Value = 4338
$ws.Range("J73").Value = 4720
$ws.Range("K73").Value = 4338
$ws.Range("L73").Value = 4720
$ws.Range("M73").Value = -3402
$ws.Range("N73").Value = -6592
$ws.Range("H102").Value = 28980.447
$ws.Range("I102").Value = 1824.1305
$ws.Range("J102").Value = 70620.13
$ws.Range("K102").Value = 1824.1305
$ws.Range("L102").Value = 70620.13
$ws.Range("M102").Value = -202.1305
$ws.Range("N102").Value = -73864.13
$ws.Range("H126").Value = 404338.6
$ws.Range("I126").Value = 3099.1667
$ws.Range("K126").Value = 9297.500100000001
$ws.Range("M126").Value = -6827.500100000001
$ws.Range("H132").Value = 3002.1296
$ws.Range("I132").Value = 2634.658
$ws.Range("K132").Value = 7903.974
$ws.Range("M132").Value = -5373.974

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 2926.6316
$ws.Range("I122").Value = 2510.6
$ws.Range("J122").Value = 3388.889
$ws.Range("K122").Value = 7531.799999999999
$ws.Range("L122").Value = 10166.667
$ws.Range("M122").Value = -5081.799999999999
$ws.Range("N122").Value = -15066.667
$ws.Range("H132").Value = 4102.8545
$ws.Range("I132").Value = 1421.3214
$ws.Range("K132").Value = 4263.9642
$ws.Range("M132").Value = -1733.9642

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 27376
$ws.Range("I96").Value = 0
$ws.Range("J96").Value = 27376
$ws.Range("K96").Value = 0
$ws.Range("L96").Value = 27376
$ws.Range("N96").Value = -30122
$ws.Range("H122").Value = 2717.6924
$ws.Range("I122").Value = 1763
$ws.Range("J122").Value = 5900
$ws.Range("K122").Value = 5289
$ws.Range("L122").Value = 17700
$ws.Range("M122").Value = -2839
$ws.Range("N122").Value = -22600
$ws.Range("M96").ClearContents()
